# Adds a new "2022-Q1" sheet (fund-level holdings) before the "总计" sheet,
# and updates the "总计" (totals) sheet with a new summary row for 2022-Q1.
#
# NOTE: sheet object references returned by Worksheets.Item(...) are NOT
# stable across structural changes (Copy/Add/rename/row insert shifts which
# underlying sheet a given reference resolves to). So every sheet handle is
# re-fetched by name immediately before it's used.

$wb = $excel.ActiveWorkbook

function Set-TextCell($wsName, $row, $col, $val) {
    # Force a numeric-looking string to be stored as TEXT (not auto-coerced
    # to a number by COM), then drop the temporary "@" number-format so the
    # cell is left with no explicit style -- matching how the sibling
    # quarter sheets store these columns.
    $ws = $wb.Worksheets.Item($wsName)
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value2 = $val
    $c.ClearFormats()
}

# ---------------------------------------------------------------------
# 1) Build the new "2022-Q1" sheet by duplicating an existing quarter
#    sheet (keeps sheetPr/margins/column styles identical) right before
#    "总计", then renaming it and overwriting its data.
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$total = $wb.Worksheets.Item("总计")
$template.Copy($total)

$q1 = $wb.Worksheets.Item("2021-Q4 (2)")
$q1.Name = "2022-Q1"

# Template had 14 data rows (2021-Q4); 2022-Q1 only has 13, drop the extra.
$q1 = $wb.Worksheets.Item("2022-Q1")
$q1.Rows.Item(15).Delete()

# Header row (row 1) B:H already reads 基金代码/基金名称/基金规模/股票总仓位/
# 仓位占比/持有市值(亿元)/仓位排名 and keeps style "2" from the template --
# nothing to change there.

$fundRows = @(
    ,@("519087", "新华优选分红混合", "9.03", "88.56", "7.12", "0.6429", 2)
    ,@("001040", "新华策略精选股票", "6.15", "93.72", "6.73", "0.4139", 2)
    ,@("519156", "新华行业轮换灵活配置混合A", "4.98", "93.77", "6.45", "0.3212", 2)
    ,@("005433", "申万菱信医药先锋股票", "2.20", "90.81", "4.83", "0.1063", 5)
    ,@("001294", "新华战略新兴产业灵活配置混合", "1.07", "93.41", "7.16", "0.0766", 2)
    ,@("011457", "新华行业龙头主题股票", "0.85", "93.55", "6.11", "0.0519", 4)
    ,@("005043", "国寿安保健康科学混合A", "0.99", "85.72", "5.19", "0.0514", 1)
    ,@("005044", "国寿安保健康科学混合C", "0.87", "85.72", "5.19", "0.0452", 1)
    ,@("501007", "汇添富中证互联网医疗主题指数（LOF）A", "0.58", "93.89", "4.74", "0.0275", 9)
    ,@("009502", "国寿安保创新医药股票A", "0.54", "81.60", "4.14", "0.0224", 3)
    ,@("501008", "汇添富中证互联网医疗主题指数（LOF）C", "0.19", "93.89", "4.74", "0.0090", 9)
    ,@("009503", "国寿安保创新医药股票C", "0.20", "81.60", "4.14", "0.0083", 3)
    ,@("519157", "新华行业轮换灵活配置混合C", "0.04", "93.77", "6.45", "0.0026", 2)
)

$r = 2
foreach ($row in $fundRows) {
    $q1 = $wb.Worksheets.Item("2022-Q1")
    # Column A is the 0-based index, numeric, style already "2" from template.
    $q1.Cells.Item($r, 1).Value2 = $r - 2
    # Columns B-G are plain text even though several look numeric.
    Set-TextCell "2022-Q1" $r 2 $row[0]
    Set-TextCell "2022-Q1" $r 3 $row[1]
    Set-TextCell "2022-Q1" $r 4 $row[2]
    Set-TextCell "2022-Q1" $r 5 $row[3]
    Set-TextCell "2022-Q1" $r 6 $row[4]
    Set-TextCell "2022-Q1" $r 7 $row[5]
    # Column H is numeric.
    $q1 = $wb.Worksheets.Item("2022-Q1")
    $q1.Cells.Item($r, 8).Value2 = $row[6]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2) Prepend a "2022-Q1" summary row to the "总计" sheet and shift the
#    existing rows (and their 0-based index column) down by one.
#
#    Rows.Insert() was tried here but it re-serialises every numeric cell
#    on the sheet at full float precision (2.45 -> 2.4500000000000002),
#    corrupting the untouched rows. Instead, every row is rewritten in
#    place top-to-bottom with its final literal values, which keeps the
#    clean decimal formatting for values that passed straight through.
# ---------------------------------------------------------------------

# Row 7 is brand new (sheet previously ended at row 6) -- give A7 the same
# style ("s=2": bold, centered, bordered) as the rest of the index column
# before anything is written into it.
$total = $wb.Worksheets.Item("总计")
$total.Range("A6").Copy()
$total = $wb.Worksheets.Item("总计")
$total.Range("A7").PasteSpecial(-4122)

$totalRows = @(
    ,@("2022-Q1", 13, 1.78)
    ,@("2021-Q4", 14, 2.45)
    ,@("2021-Q3", 21, 8.49)
    ,@("2021-Q2", 22, 7.62)
    ,@("2021-Q1", 12, 3.94)
    ,@("2020-Q4", 7, 1.83)
)

$row = 2
foreach ($entry in $totalRows) {
    $total = $wb.Worksheets.Item("总计")
    $total.Cells.Item($row, 1).Value2 = $row - 2
    $total.Cells.Item($row, 2).Value2 = $entry[0]
    $total.Cells.Item($row, 3).Value2 = $entry[1]
    $total.Cells.Item($row, 4).Value2 = $entry[2]
    $row = $row + 1
}
